# Fruta / hortaliza, semanal
# Insert a new weekly record at row 391 (pushing the existing rows 391..473
# down to 392..474) on the single worksheet of the workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 391; this shifts rows 391-473 down to 392-474
# and automatically extends the used range / dimension to A1:R474.
$ws.Rows.Item(391).Insert()

# Populate the newly inserted row 391 with the new weekly observation.
$ws.Cells.Item(391, 1).Value = 3
$ws.Cells.Item(391, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(391, 3).Value = "Coquimbo"
$ws.Cells.Item(391, 4).Value2 = 44889
$ws.Cells.Item(391, 5).Value = 5
$ws.Cells.Item(391, 6).Value = 100112013
$ws.Cells.Item(391, 7).Value = "Alcachofa"
$ws.Cells.Item(391, 8).Value = "Española"
$ws.Cells.Item(391, 9).Value = "Primera"
$ws.Cells.Item(391, 10).Value = 25000
$ws.Cells.Item(391, 11).Value = 200
$ws.Cells.Item(391, 12).Value = 220
$ws.Cells.Item(391, 13).Value = 210
$ws.Cells.Item(391, 14).Value = "$/unidad"
$ws.Cells.Item(391, 15).Value = "Llay Llay"
$ws.Cells.Item(391, 16).Value = 210
$ws.Cells.Item(391, 17).Value = 1
$ws.Cells.Item(391, 18).Value = "Hortaliza"
